$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Added paddle shifter PCB and the relays.
#
# New wires (rows 73-78) describing the paddle-shifter gearbox pneumatics
# board (pressure sensor, switch/up/down solenoids, pump, neutral sensor)
# plus a note about the dash wiring (row 81).
#
# The cell values below are written in the same order the original author
# entered them so brand-new shared-string entries land at the same indices
# as the target workbook (existing strings are simply re-used / deduped).
# ---------------------------------------------------------------------------

# -- Column A labels for the six new wires -----------------------------
$ws.Cells.Item(73,1).Value2 = "ShifterPressure"
$ws.Cells.Item(74,1).Value2 = "SwitchSolenoid"
$ws.Cells.Item(75,1).Value2 = "UpSolenoid"
$ws.Cells.Item(76,1).Value2 = "DownSolenoid"
$ws.Cells.Item(77,1).Value2 = "ShifterPump"
$ws.Cells.Item(78,1).Value2 = "NeutralSensor"

# -- Neutral sensor row's connection / description ----------------------
$ws.Cells.Item(78,4).Value2 = "Neutral Sensor"
$ws.Cells.Item(78,12).Value2 = "Signal of the neutral sensor on the engine"

# -- Shared "Connection3" for the pneumatics wires -----------------------
$ws.Cells.Item(73,4).Value2 = "Gearbox  pneumatics"

# -- Descriptions (column L) for the solenoid / pump / pressure wires ---
$ws.Cells.Item(76,12).Value2 = "Output to solenoid on down shift valve"
$ws.Cells.Item(75,12).Value2 = "Output to solenoid on up shift valve"
$ws.Cells.Item(74,12).Value2 = "Output to switch solenoid used to start pump from atmosphere"
$ws.Cells.Item(77,12).Value2 = "High current pump output"
$ws.Cells.Item(73,12).Value2 = "Analog input from pneumatic pressure sensor"

# -- Current ratings (column M) ------------------------------------------
$ws.Cells.Item(74,13).Value2 = "0.5A"
$ws.Cells.Item(77,13).Value2 = "6A"

# -- Dash note on row 81 --------------------------------------------------
$ws.Cells.Item(81,1).Value2 = "Wires in the dash area are still very much WIP."

# ---------------------------------------------------------------------------
# Remaining columns for rows 73-78 (re-use existing shared strings / plain
# numbers, so no new entries are minted).
# ---------------------------------------------------------------------------

# Connection1 (device) - all new wires belong to the paddle shifters device
$ws.Cells.Item(73,2).Value2 = "Paddle shifters"
$ws.Cells.Item(74,2).Value2 = "Paddle shifters"
$ws.Cells.Item(75,2).Value2 = "Paddle shifters"
$ws.Cells.Item(76,2).Value2 = "Paddle shifters"
$ws.Cells.Item(77,2).Value2 = "Paddle shifters"
$ws.Cells.Item(78,2).Value2 = "Paddle shifters"

# Pin numbers (column C)
$ws.Cells.Item(73,3).Value2 = 14
$ws.Cells.Item(74,3).Value2 = 5
$ws.Cells.Item(75,3).Value2 = 3
$ws.Cells.Item(76,3).Value2 = 4
$ws.Cells.Item(77,3).Value2 = 6
$ws.Cells.Item(78,3).Value2 = 13

# Connection3 (column D) for the solenoid/pump/pressure wires
$ws.Cells.Item(74,4).Value2 = "Gearbox  pneumatics"
$ws.Cells.Item(75,4).Value2 = "Gearbox  pneumatics"
$ws.Cells.Item(76,4).Value2 = "Gearbox  pneumatics"
$ws.Cells.Item(77,4).Value2 = "Gearbox  pneumatics"

# Connection4 (column E) - all TBD
$ws.Cells.Item(73,5).Value2 = "TBD"
$ws.Cells.Item(74,5).Value2 = "TBD"
$ws.Cells.Item(75,5).Value2 = "TBD"
$ws.Cells.Item(76,5).Value2 = "TBD"
$ws.Cells.Item(77,5).Value2 = "TBD"
$ws.Cells.Item(78,5).Value2 = "TBD"

# Current (column M) for remaining rows
$ws.Cells.Item(73,13).Value2 = "Signal"
$ws.Cells.Item(75,13).Value2 = "0.5A"
$ws.Cells.Item(76,13).Value2 = "0.5A"
$ws.Cells.Item(78,13).Value2 = "Signal"

# Can be carried over? (column O)
$ws.Cells.Item(73,15).Value2 = "No, does not exist"
$ws.Cells.Item(74,15).Value2 = "No, does not exist"
$ws.Cells.Item(75,15).Value2 = "No, does not exist"
$ws.Cells.Item(76,15).Value2 = "No, does not exist"
$ws.Cells.Item(77,15).Value2 = "No, does not exist"
$ws.Cells.Item(78,15).Value2 = "No, does not exist"

# The first pin entry (C73) carries a distinct number-format style in the
# source workbook; re-apply "General" so the cell gets its own style entry.
$ws.Cells.Item(73,3).NumberFormat = "General"

# Leave the selection on the newly added note, matching the author's final
# cursor position.
$ws.Range("A81").Select()
